# Apply "difference instance" adjustments to schedule_6.xlsx
# - Column B (Flight No.) becomes a numeric row index (same as column A) instead of a flight-code string
# - Columns C and D (Arrival/Departure Time) get new time values
# - Column E (Location) gets new location codes
# - Column F (Type) gets new values (in some rows)
# Columns A, G and H are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  B=1;  C=0.03194444444444444; D=0.05694444444444444; E='D37'; F=3 },
    @{ Row=3;  B=2;  C=0.009722222222222222; D=0.03055555555555555; E='B8'; F=1 },
    @{ Row=4;  B=3;  C=0.02361111111111111; D=0.04861111111111111; E='E11'; F=3 },
    @{ Row=5;  B=4;  C=0.03333333333333333; D=0.05625; E='A20'; F=3 },
    @{ Row=6;  B=5;  C=0.0375; D=0.05972222222222223; E='A7'; F=2 },
    @{ Row=7;  B=6;  C=0.02569444444444444; D=0.04930555555555555; E='F37'; F=2 },
    @{ Row=8;  B=7;  C=0.02569444444444444; D=0.04930555555555555; E='F41'; F=2 },
    @{ Row=9;  B=8;  C=0.03888888888888889; D=0.0625; E='C20'; F=1 },
    @{ Row=10; B=9;  C=0.03958333333333333; D=0.06597222222222222; E='B9'; F=3 },
    @{ Row=11; B=10; C=0.009722222222222222; D=0.03194444444444444; E='F34'; F=1 },
    @{ Row=12; B=11; C=0.04097222222222222; D=0.0625; E='E21'; F=3 },
    @{ Row=13; B=12; C=0.02986111111111111; D=0.05694444444444444; E='D32'; F=3 },
    @{ Row=14; B=13; C=0.02013888888888889; D=0.04375; E='A15'; F=3 },
    @{ Row=15; B=14; C=0.025; D=0.05138888888888889; E='E7'; F=2 },
    @{ Row=16; B=15; C=0.01944444444444444; D=0.04097222222222222; E='F35'; F=1 },
    @{ Row=17; B=16; C=0.02152777777777778; D=0.04375; E='D41'; F=2 },
    @{ Row=18; B=17; C=0.025; D=0.04930555555555555; E='C18'; F=2 },
    @{ Row=19; B=18; C=0.02222222222222222; D=0.04305555555555556; E='C22'; F=2 },
    @{ Row=20; B=19; C=0.04027777777777778; D=0.06388888888888888; E='A2'; F=2 },
    @{ Row=21; B=20; C=0.02569444444444444; D=0.04861111111111111; E='F52'; F=3 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
